$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.529.88'
$ws.Range('D3').Value = '1.623.12'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''211.72'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').Value = '''0.525'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '''23.13'
$ws.Range('E8').Value = '  -1.26%  '
$ws.Range('E9').Value = '  +1.49%  '
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('D11').Value = '''0.0879'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').Value = '1.853.42'
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('D13').Value = '1.619.13'
$ws.Range('E13').Value = '  -1.52%  '
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('E15').Value = '  -2.00%  '
$ws.Range('D16').Value = '''65.17'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '27.491.60'
$ws.Range('E17').Value = '  -0.63%  '
$ws.Range('D18').Value = '''229.90'
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('D19').Value = '0.0₃0718'
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('E20').Value = '  -1.42%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  +3.32%  '
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('E24').Value = '  +5.06%  '
$ws.Range('E26').Value = '  -1.07%  '
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('D33').Value = '1.469.43'
$ws.Range('E33').Value = '  +1.74%  '
$ws.Range('D34').Value = '''3.05'
$ws.Range('E34').Value = '  -2.58%  '
$ws.Range('D35').Value = '''1.54'
$ws.Range('E35').Value = '  -2.73%  '
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('D37').Value = '''0.935'
$ws.Range('E37').Value = '  +5.10%  '
$ws.Range('D38').Value = '''0.871'
$ws.Range('E38').Value = '  -1.05%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.0166'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = '''0.555'
$ws.Range('E40').Value = '  -2.69%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  -1.94%  '
$ws.Range('D43').Value = '''67.18'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D45').Value = '''2.20'
$ws.Range('E45').Value = '  -2.19%  '
$ws.Range('D46').Value = '''5.29'
$ws.Range('E46').Value = '  -6.07%  '
$ws.Range('D47').Value = '1.763.27'
$ws.Range('E47').Value = '  -1.33%  '
$ws.Range('E48').Value = '  +1.41%  '
$ws.Range('D49').Value = '''87.30'
$ws.Range('E49').Value = '  +2.00%  '
$ws.Range('E50').Value = '  -1.31%  '
$ws.Range('D51').Value = '''0.0998'
$ws.Range('E51').Value = '  +0.92%  '

Write-Host "Applied all cryptos updates"
